# Update the monthly sales/revenue/expenses table:
#  - rename months to "Mon-2023" style labels and extend through December
#  - refresh sales/revenue numbers for the (now) 12 months
#  - turn the "expenses" column into a B-C formula
#  - move the cell selection to F3 (as last left by the author)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("month") must stay plain text - format it as Text *before* writing
# the "Jan-2023" style values so Excel doesn't auto-coerce them into dates.
$ws.Range("A1:A13").NumberFormat = "@"

$months  = @("Jan-2023","Feb-2023","Mar-2023","Apr-2023","May-2023","June-2023","July-2023","Aug-2023","Sep-2023","Oct-2023","Nov-2023","Dec-2023")
$sales   = @(98,42,28,40,91,60,90,37,36,24,60,32)
$revenue = @(50,18,27,30,54,54,3,30,20,10,20,25)

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = $sales[$i]
    $ws.Cells.Item($row, 3).Value = $revenue[$i]
}

# "expenses" column: D2 gets its own formula, D3:D13 are written together so
# they become one shared-formula group (mirrors a fill-down in the UI).
$ws.Range("D2").Formula = "=B2-C2"
$ws.Range("D3:D13").Formula = "=B3-C3"

# Leave the selection where the author left it.
[void]$ws.Range("F3").Select()
